$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date serial for every
# data row (2..285). The commit bumped that date by one day (46082 ->
# 46083, i.e. 2026-03-01 -> 2026-03-02) across the entire column.
$ws.Range("C2:C285").Value = 46083
